# Remote Indigenous Housing - state views.
#
# Updates the "Description" sheet of the housing_remote_indigenous workbook:
#  - Shortens the "Other Benchmarks" paragraph (drops the NT sentence, which
#    moves into its own per-state breakdown below).
#  - Re-purposes the row that used to hold the Qld paragraph for "Notes".
#  - Adds a new per-jurisdiction breakdown (ACT / Vic / Tas / NT / Qld /
#    Australia) below the existing notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# --- Row 5: "Other Benchmarks" body shortened (drop the NT sentence) -------
$ws.Cells.Item(5, 1).Value = "Other Benchmarks"
$ws.Cells.Item(5, 2).Value = "Under NPARIH, states and the Northern Territory have a 20 per cent Indigenous employment target for capital works.  Over the life of NPARIH, all jurisdictions have met or exceeded this target."
$ws.Cells.Item(5, 2).WrapText = $true
$ws.Rows.Item(5).RowHeight = 35.05

# --- Row 6: now "Notes" / NPARIH notional targets ---------------------------
$ws.Cells.Item(6, 1).Value = "Notes"
$ws.Cells.Item(6, 2).Value = "NPARIH notional targets extend to 2014 for refurbishments and to 2018 for new build houses."
$ws.Cells.Item(6, 2).WrapText = $true
$ws.Rows.Item(6).RowHeight = 23.95

# --- Row 7: Victoria not included note --------------------------------------
$ws.Cells.Item(7, 1).Value = ""
$ws.Cells.Item(7, 2).Value = "Victoria is not included in these analyses as no Commonwealth funding was provided to Victoria under the NPARIH for the construction of new houses or refurbishment of existing stock."
$ws.Cells.Item(7, 2).WrapText = $true
$ws.Rows.Item(7).RowHeight = 35.05

# --- Row 8: Victoria and Tasmania exited note -------------------------------
$ws.Cells.Item(8, 1).Value = ""
$ws.Cells.Item(8, 2).Value = "Victoria and Tasmania are not currently a part of NPARIH as they exited in 2014."
$ws.Cells.Item(8, 2).WrapText = $true
$ws.Rows.Item(8).RowHeight = 23.85

# --- Row 9: ACT --------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = "ACT"
$ws.Cells.Item(9, 2).Value = "ACT does not participate in this agreement"
$ws.Rows.Item(9).RowHeight = 12.8

# --- Row 10: Vic --------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = "Vic"
$ws.Cells.Item(10, 2).Value = "Victoria existed this agreement in 2014."
$ws.Rows.Item(10).RowHeight = 12.8

# --- Row 11: Tas --------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = "Tas"
$ws.Cells.Item(11, 2).Value = "Tasmania exited this agreement in 2014."
$ws.Rows.Item(11).RowHeight = 12.8

# --- Row 12: NT ----------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = "NT"
$ws.Cells.Item(12, 2).Value = "The Northern Territory Government has consistently delivered against a 40 per cent Indigenous employment target for property management and a 40 per cent Indigenous employment target for tenancy management, as outlined in the Northern Territory’s 2014-16 Implementation Plan to the NPARIH."
$ws.Cells.Item(12, 2).WrapText = $true
$ws.Rows.Item(12).RowHeight = 46.45

# --- Row 13: Qld -----------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Qld"
$ws.Cells.Item(13, 2).Value = "Queensland has consistently exceeded the 20 per cent Indigenous employment target for capital works, with over 85 per cent of all available jobs in construction, repairs and maintenance during 2014-15 undertaken by Aboriginal and Torres Strait Islander workers."
$ws.Cells.Item(13, 2).WrapText = $true
$ws.Rows.Item(13).RowHeight = 46.45

# --- Row 14: Australia -------------------------------------------------------
$ws.Cells.Item(14, 1).Value = "Australia"
$ws.Cells.Item(14, 2).Value = "In addition, the Northern Territory Government has consistently delivered against a 40 per cent Indigenous employment target for property management and a 40 per cent Indigenous employment target for tenancy management, as outlined in the Northern Territory’s 2014-16 Implementation Plan to the NPARIH."
$ws.Cells.Item(14, 2).WrapText = $true
$ws.Rows.Item(14).RowHeight = 57.7

# --- Row 15: repeats the Qld paragraph --------------------------------------
$ws.Cells.Item(15, 2).Value = "Queensland has consistently exceeded the 20 per cent Indigenous employment target for capital works, with over 85 per cent of all available jobs in construction, repairs and maintenance during 2014-15 undertaken by Aboriginal and Torres Strait Islander workers."
$ws.Cells.Item(15, 2).WrapText = $true
$ws.Rows.Item(15).RowHeight = 46.45

# --- Column widths (slightly narrower after the edit) ------------------------
$ws.Columns.Item(1).ColumnWidth = 17.43
$ws.Columns.Item(2).ColumnWidth = 62.29

# --- Scroll the sheet down a little so row 4 sits at the top -----------------
$ws.Activate()
$ws.Range("A6").Select()

# --- Tab-bar/scrollbar split ratio tweak seen on the workbook window ---------
$excel.ActiveWindow.TabRatio = 0.5
